$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every cell (even numeric-looking ones like "200") as
# literal TEXT, not as a number. A plain `.Value = "200"` assignment would
# be auto-coerced to a real number by Excel, and a plain `.Value = ""`
# assignment clears/removes the cell instead of leaving an empty text
# cell behind. To avoid both pitfalls we stage each value as a text
# formula in a scratch cell, then Cut it (a "move", not a copy) into its
# final destination - this drops the formula and leaves a plain literal
# value in place, without forcing a number interpretation.
function Set-TextValue {
    param($cell, [string]$text)
    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Cut($cell)
}

Set-TextValue $ws.Range("A5") ""
Set-TextValue $ws.Range("B5") "حسن "
Set-TextValue $ws.Range("C5") "200"
Set-TextValue $ws.Range("D5") "النصر"
Set-TextValue $ws.Range("E5") "الرحلة 3"
Set-TextValue $ws.Range("F5") "C4"
Set-TextValue $ws.Range("G5") "NRC"
Set-TextValue $ws.Range("H5") "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٢٩:٣٠ م"
